{"js": "// Insert two new paragraphs right after the paragraph that ends the\n// \"Model Services: ...\" bullet and before the \"Outputs. Connectors ...\"\n// bullet: a blank paragraph, then the new \"Interactions: Services. ...\"\n// paragraph (matching the existing document's empty-line-between-bullets\n// pattern).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst targetText =\n  \"Model Services: Monads parsed DOM interactions services (contexts). Render / update DOM.\";\nconst newText =\n  \"Interactions: Services. Browse DOM. Apply selectors / Browse available transforms (Monads / HATEOAS).\";\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === targetText) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Model Services: ...' paragraph.\");\n}\n\n// Insert an empty paragraph directly after the target, then the new\n// text paragraph directly after that empty one (so the final order is:\n// Model Services..., <blank>, Interactions: Services..., <blank>, Outputs...).\nconst blankPara = target.insertParagraph(\"\", \"After\");\nblankPara.insertParagraph(newText, \"After\");\n\nawait context.sync();\n", "ps1": "# Insert two new paragraphs right after the \"Model Services: ...\" bullet\n# and before the \"Outputs. Connectors ...\" bullet: a blank paragraph,\n# then the new \"Interactions: Services. ...\" paragraph (matching the\n# existing document's empty-line-between-bullets pattern).\n\n$d = $word.ActiveDocument\n\n$targetText = \"Model Services: Monads parsed DOM interactions services (contexts). Render / update DOM.\"\n$newText = \"Interactions: Services. Browse DOM. Apply selectors / Browse available transforms (Monads / HATEOAS).\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq $targetText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Model Services: ...' paragraph.\"\n}\n\n# First insert: a blank paragraph right after the target.\n$target.Range.InsertParagraphAfter()\n$blank = $target.Next()\n\n# Second insert: the new text paragraph right after the blank one.\n$blank.Range.InsertParagraphAfter()\n$textPara = $blank.Next()\n$textPara.Range.Text = $newText\n"}
